# Actualización automática 2025-07-09 12:30:08
# Updates commission/sales figures for HIDALGO HIDALGO PEDRO GUSTAVO
# across the three report sheets (VENTAS POR GRUPO, VENTA MENSUAL,
# CUMPLIMIENTO MENSUAL) to reflect newly recorded PORCELANATO sales.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M16").Value = 6978.81
$wsGrupo.Range("M21").Value = 2156.54
$wsGrupo.Range("M22").Value = "7 de 20"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F16").Value = 6978.81
$wsMensual.Range("F21").Value = 2156.54
$wsMensual.Range("F22").Value = 25684.6

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D16").Value = 19855.88
$wsCumplimiento.Range("E16").Value = 24410.36
$wsCumplimiento.Range("F16").Value = 0.4485558294537779
$wsCumplimiento.Range("D19").Value = 25684.6
$wsCumplimiento.Range("E19").Value = 39693.39762291769
$wsCumplimiento.Range("F19").Value = 0.3928630568978528
